# Update cryptocurrency price/volume data (Price column D, Volume(1h) column E)
# Leading apostrophe forces text interpretation (some prices parse as numbers,
# e.g. "493.42"); Style reset to "Normal" afterward keeps cell formatting
# identical to the source (no numFmt/style churn from the text coercion).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.646.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.491.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'493.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.70%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'152.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +7.00%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.88%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.504.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.54%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'5.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.67%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.27%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.73%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.77%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.922.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.32%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'56.766.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.46%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.74%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.501.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.90%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.77%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.23%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'320.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.57%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'58.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.62%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.54%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.601.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.87%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.72%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.55%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.11%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'151.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.06%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'18.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.20%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.50%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.98%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.42%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.09%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.00%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'34.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.85%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0566"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.71%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.617"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.36%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'4.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.11%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'268.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0929"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.47%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0230"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.40%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'10.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.91%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'17.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.07%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.890.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.43%  "
$ws.Range("E51").Style = "Normal"
